# dsa heaps and math
# Adds two new rows (355. Design Twitter / Heaps, 66. Plus One / Math)
# to the LeetCode tracking table on Sheet1, expanding the Excel table,
# dimension, hyperlinks and shared strings accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the existing table (Table2) from A1:E85 to A1:E87 so the new rows
# become part of the structured table / autofilter range.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E87"))

# ---- Row 86: 355. Design Twitter ----------------------------------------
$ws.Range("A86").Value = "355. Design Twitter"

$ws.Range("B86").Value = "Medium"
$ws.Range("B86").Interior.Color = 49407

$ws.Range("C86").Value = "Heaps"

$ws.Range("D86").Value = "Review. Similar to Merge K Sorted List. 1. Follow/unfollow functions HashMap of UserIds, maps to a HashSet of followeeIds. 2. postTweet is a HashMap of userIds which maps to a list of [count, TweetIds]. 3. Put the time of tweets in a Max Heap for getNewsFeed."

$ws.Range("E86").Value = "https://leetcode.com/problems/design-twitter/solutions/82825/java-oo-design-with-most-efficient-function-getnewsfeed/ "
$ws.Hyperlinks.Add($ws.Range("E86"), "https://leetcode.com/problems/design-twitter/solutions/82825/java-oo-design-with-most-efficient-function-getnewsfeed/") | Out-Null
$ws.Range("E86").Style = "Hyperlink"

# ---- Row 87: 66. Plus One ------------------------------------------------
$ws.Range("A87").Value = "66. Plus One"

$ws.Range("B87").Value = "Easy"
$ws.Range("B87").Interior.Color = 5287936

$ws.Range("C87").Value = "Math"

$ws.Range("D87").Value = "Remember the carry when reaching 10. The crux is to manage 2 potential result arrays, newDigits and the input digits array in place. If there is a carry at the end, we add it to newDigits[0] and then return it, else we just return the input array which we computed in place."

$ws.Range("E87").Value = "https://leetcode.com/problems/plus-one/solutions/2706861/java-fastest-0ms-runtime-easy-and-elegant-solution/ "
$ws.Hyperlinks.Add($ws.Range("E87"), "https://leetcode.com/problems/plus-one/solutions/2706861/java-fastest-0ms-runtime-easy-and-elegant-solution/") | Out-Null
$ws.Range("E87").Style = "Hyperlink"

# ---- Update the view so the newly added rows are in focus ---------------
$ws.Range("E93").Select() | Out-Null
